$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            try {
                if ($sh.TextFrame.TextRange.Text -ne $newText) {
                    $sh.TextFrame.TextRange.Text = $newText
                }
            } catch {
            }
        }
    }
}

# Update the Date placeholder cached text on the slide master
Update-DatePlaceholder $p.SlideMaster.Shapes "19/07/2018"

# Update the Date placeholder cached text on every slide layout
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholder $lay.Shapes "19/07/2018"
}
